$wb = $excel.ActiveWorkbook

# Rename the second sheet ("connected_areaType" -> "Catchment_LanduseMix")
$wsLanduse = $wb.Worksheets.Item("connected_areaType")
$wsLanduse.Name = "Catchment_LanduseMix"

$wsSite = $wb.Worksheets.Item("site_data")

# Update manually-entered data values on the site_data sheet
$wsSite.Range("C6").Value = 0.04
$wsSite.Range("C14").Value = 5.2619999999999996

# Correct the spelling of the header on the Catchment_LanduseMix sheet
$wsLanduse.Range("D1").Value = "separate_sewer_percent"

# Leave the selection on site_data parked at D31
$wsSite.Range("D31").Select()

# Switch to the Catchment_LanduseMix sheet and select D2, making it the active tab
$wsLanduse.Activate()
$wsLanduse.Range("D2").Select()
